$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 337
$ws.Range("I9").Value = 179.2
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 179.2
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = -10.19999999999999
$ws.Range("N9").Value = -938
$ws.Range("H19").Value = 1638.762
$ws.Range("I19").Value = 1860
$ws.Range("J19").Value = 1528.1428
$ws.Range("K19").Value = 1860
$ws.Range("L19").Value = 1528.1428
$ws.Range("M19").Value = -1685
$ws.Range("N19").Value = -1878.1428
$ws.Range("H38").Value = 563.6667
$ws.Range("I38").Value = 563.6667
$ws.Range("K38").Value = 1691.0001
$ws.Range("M38").Value = -1319.0001
$ws.Range("H42").Value = 753.9
$ws.Range("J42").Value = 1042.2858
$ws.Range("L42").Value = 3126.8574
$ws.Range("N42").Value = -3586.8574
$ws.Range("H43").Value = 6162
$ws.Range("I43").Value = 2250
$ws.Range("J43").Value = 10074
$ws.Range("K43").Value = 2250
$ws.Range("L43").Value = 10074
$ws.Range("M43").Value = -2181
$ws.Range("N43").Value = -10212
$ws.Range("H70").Value = 51000
$ws.Range("I70").Value = 100000
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 300000
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -299730
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 51000
$ws.Range("I73").Value = 100000
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 300000
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -299064
$ws.Range("N73").Value = -7872
$ws.Range("H80").Value = 580.5833
$ws.Range("I80").Value = 601.55554
$ws.Range("J80").Value = 517.6667
$ws.Range("K80").Value = 1804.66662
$ws.Range("L80").Value = 1553.0001
$ws.Range("M80").Value = -806.66662
$ws.Range("N80").Value = -3549.0001
$ws.Range("H83").Value = 580.5833
$ws.Range("I83").Value = 601.55554
$ws.Range("J83").Value = 517.6667
$ws.Range("K83").Value = 5413.99986
$ws.Range("L83").Value = 4659.0003
$ws.Range("M83").Value = -421.9998599999999
$ws.Range("N83").Value = -14643.0003
$ws.Range("H86").Value = 3782
$ws.Range("I86").Value = 3638.4
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 3638.4
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -2515.4
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 3782
$ws.Range("I89").Value = 3638.4
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 18192
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -12576
$ws.Range("N89").Value = -33732
$ws.Range("H100").Value = 1618.6666
$ws.Range("I100").Value = 1321.88
$ws.Range("J100").Value = 2546.125
$ws.Range("K100").Value = 1321.88
$ws.Range("L100").Value = 2546.125
$ws.Range("M100").Value = -780.8800000000001
$ws.Range("N100").Value = -3628.125
$ws.Range("H105").Value = 40671
$ws.Range("J105").Value = 40671
$ws.Range("L105").Value = 40671
$ws.Range("N105").Value = -47659
$ws.Range("H132").Value = 3096.5715
$ws.Range("I132").Value = 1400.7812
$ws.Range("K132").Value = 4202.3436
$ws.Range("M132").Value = -1672.3436
$ws.Range("H135").Value = 834.0714
$ws.Range("I135").Value = 689.8333
$ws.Range("K135").Value = 6208.4997
$ws.Range("M135").Value = -3673.4997
$ws.Range("H137").Value = 3337.5862
$ws.Range("I137").Value = 2652.3845
$ws.Range("J137").Value = 3894.3125
$ws.Range("K137").Value = 7957.1535
$ws.Range("L137").Value = 11682.9375
$ws.Range("M137").Value = -5407.1535
$ws.Range("N137").Value = -16782.9375
$ws.Range("H138").Value = 2099.4285
$ws.Range("J138").Value = 4129.143
$ws.Range("L138").Value = 12387.429
$ws.Range("N138").Value = -22667.429
$ws.Range("H141").Value = 129528
$ws.Range("I141").Value = 129528
$ws.Range("K141").Value = 388584
$ws.Range("M141").Value = -383404

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 110712
$ws.Range("J7").Value = 110712
$ws.Range("L7").Value = 110712
$ws.Range("N7").Value = -110940
$ws.Range("H32").Value = 4519847
$ws.Range("I32").Value = 763697
$ws.Range("J32").Value = 27891446
$ws.Range("K32").Value = 763697
$ws.Range("L32").Value = 27891446
$ws.Range("M32").Value = -763410
$ws.Range("N32").Value = -27892020
$ws.Range("H63").Value = 1897.6154
$ws.Range("I63").Value = 1805.8334
$ws.Range("J63").Value = 2999
$ws.Range("K63").Value = 1805.8334
$ws.Range("L63").Value = 2999
$ws.Range("M63").Value = -1119.8334
$ws.Range("N63").Value = -4371
$ws.Range("H66").Value = 1897.6154
$ws.Range("I66").Value = 1805.8334
$ws.Range("J66").Value = 2999
$ws.Range("K66").Value = 9029.166999999999
$ws.Range("L66").Value = 14995
$ws.Range("M66").Value = -5597.166999999999
$ws.Range("N66").Value = -21859
$ws.Range("H74").Value = 2688.923
$ws.Range("I74").Value = 1633
$ws.Range("J74").Value = 5064.75
$ws.Range("K74").Value = 1633
$ws.Range("L74").Value = 5064.75
$ws.Range("M74").Value = -759
$ws.Range("N74").Value = -6812.75
$ws.Range("H77").Value = 2688.923
$ws.Range("I77").Value = 1633
$ws.Range("J77").Value = 5064.75
$ws.Range("K77").Value = 8165
$ws.Range("L77").Value = 25323.75
$ws.Range("M77").Value = -3797
$ws.Range("N77").Value = -34059.75
$ws.Range("H97").Value = 25720.555
$ws.Range("I97").Value = 31426.428
$ws.Range("K97").Value = 31426.428
$ws.Range("M97").Value = -30930.428
$ws.Range("H102").Value = 1392.8
$ws.Range("I102").Value = 1339.25
$ws.Range("K102").Value = 1339.25
$ws.Range("M102").Value = 282.75
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H118").Value = 60000
$ws.Range("J118").Value = 60000
$ws.Range("L118").Value = 60000
$ws.Range("N118").Value = -63314
$ws.Range("H122").Value = 1999.6666
$ws.Range("I122").Value = 1799.5
$ws.Range("K122").Value = 5398.5
$ws.Range("M122").Value = -2948.5
$ws.Range("H132").Value = 3660.8147
$ws.Range("I132").Value = 3660.8147
$ws.Range("K132").Value = 10982.4441
$ws.Range("M132").Value = -8452.444100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 20889
$ws.Range("J74").Value = 20889
$ws.Range("L74").Value = 20889
$ws.Range("N74").Value = -22761
$ws.Range("H77").Value = 20889
$ws.Range("J77").Value = 20889
$ws.Range("L77").Value = 62667
$ws.Range("N77").Value = -72027
$ws.Range("H86").Value = 1962.9445
$ws.Range("J86").Value = 2238.8
$ws.Range("L86").Value = 2238.8
$ws.Range("N86").Value = -4484.8
$ws.Range("H89").Value = 1962.9445
$ws.Range("J89").Value = 2238.8
$ws.Range("L89").Value = 11194
$ws.Range("N89").Value = -22426
$ws.Range("H94").Value = 2068.1304
$ws.Range("I94").Value = 1550.9048
$ws.Range("K94").Value = 1550.9048
$ws.Range("M94").Value = -1099.9048
$ws.Range("H99").Value = 5013.5
$ws.Range("I99").Value = 4750
$ws.Range("J99").Value = 5277
$ws.Range("K99").Value = 4750
$ws.Range("L99").Value = 5277
$ws.Range("M99").Value = -3252
$ws.Range("N99").Value = -8273
$ws.Range("H107").Value = 4122.5415
$ws.Range("I107").Value = 3925.1875
$ws.Range("K107").Value = 3925.1875
$ws.Range("M107").Value = -2005.1875
$ws.Range("H108").Value = 98765
$ws.Range("J108").Value = 98765
$ws.Range("L108").Value = 98765
$ws.Range("N108").Value = -106445

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 188.3125
$ws.Range("I7").Value = 119.63636
$ws.Range("K7").Value = 119.63636
$ws.Range("M7").Value = -6.636359999999996
$ws.Range("H22").Value = 697.1429000000001
$ws.Range("I22").Value = 296.5
$ws.Range("J22").Value = 1231.3334
$ws.Range("K22").Value = 296.5
$ws.Range("L22").Value = 1231.3334
$ws.Range("M22").Value = 53.5
$ws.Range("N22").Value = -1931.3334
$ws.Range("H31").Value = 2162.6985
$ws.Range("I31").Value = 1464.5366
$ws.Range("J31").Value = 3463.818
$ws.Range("K31").Value = 1464.5366
$ws.Range("L31").Value = 3463.818
$ws.Range("M31").Value = -1169.5366
$ws.Range("N31").Value = -4053.818
$ws.Range("H34").Value = 2162.6985
$ws.Range("I34").Value = 1464.5366
$ws.Range("J34").Value = 3463.818
$ws.Range("K34").Value = 1464.5366
$ws.Range("L34").Value = 3463.818
$ws.Range("M34").Value = -1262.5366
$ws.Range("N34").Value = -3867.818
$ws.Range("H103").Value = 16404.5
$ws.Range("I103").Value = 16404.5
$ws.Range("K103").Value = 16404.5
$ws.Range("M103").Value = -15232.5
$ws.Range("H107").Value = 2322.6667
$ws.Range("I107").Value = 2194.4546
$ws.Range("J107").Value = 2675.25
$ws.Range("K107").Value = 2194.4546
$ws.Range("L107").Value = 2675.25
$ws.Range("M107").Value = -274.4546
$ws.Range("N107").Value = -6515.25
$ws.Range("H132").Value = 1114.5217
$ws.Range("I132").Value = 839.7619
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 2519.2857
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = 10.71430000000009
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 830.7931
$ws.Range("I134").Value = 703.4400000000001
$ws.Range("K134").Value = 2110.32
$ws.Range("M134").Value = 424.6799999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 33593.965
$ws.Range("J11").Value = 204.25
$ws.Range("L11").Value = 612.75
$ws.Range("N11").Value = -892.75
$ws.Range("H12").Value = 705.44446
$ws.Range("I12").Value = 596.8
$ws.Range("K12").Value = 1790.4
$ws.Range("M12").Value = -1617.4
$ws.Range("H14").Value = 290.55554
$ws.Range("I14").Value = 290.55554
$ws.Range("K14").Value = 871.66662
$ws.Range("M14").Value = -698.66662
$ws.Range("H38").Value = 66.45
$ws.Range("I38").Value = 39.066666
$ws.Range("J38").Value = 148.6
$ws.Range("K38").Value = 117.199998
$ws.Range("L38").Value = 445.8
$ws.Range("M38").Value = 229.800002
$ws.Range("N38").Value = -1139.8
$ws.Range("H68").Value = 1668181.9
$ws.Range("I68").Value = 1755.1111
$ws.Range("K68").Value = 5265.3333
$ws.Range("M68").Value = -4454.3333
$ws.Range("H71").Value = 1668181.9
$ws.Range("I71").Value = 1755.1111
$ws.Range("K71").Value = 15795.9999
$ws.Range("M71").Value = -11739.9999
$ws.Range("H99").Value = 11635.143
$ws.Range("I99").Value = 2361.5
$ws.Range("J99").Value = 24000
$ws.Range("K99").Value = 7084.5
$ws.Range("L99").Value = 72000
$ws.Range("M99").Value = -4838.5
$ws.Range("N99").Value = -76492
$ws.Range("H107").Value = 982.6842
$ws.Range("I107").Value = 591.6667
$ws.Range("J107").Value = 1056
$ws.Range("K107").Value = 1775.0001
$ws.Range("L107").Value = 3168
$ws.Range("M107").Value = 144.9999
$ws.Range("N107").Value = -7008
$ws.Range("H129").Value = 68074.336
$ws.Range("I129").Value = 154267
$ws.Range("J129").Value = 2162.2942
$ws.Range("K129").Value = 462801
$ws.Range("L129").Value = 6486.882599999999
$ws.Range("M129").Value = -457801
$ws.Range("N129").Value = -16486.8826
$ws.Range("H131").Value = 17315772
$ws.Range("I131").Value = 35793684
$ws.Range("J131").Value = 69717.87
$ws.Range("K131").Value = 107381052
$ws.Range("L131").Value = 209153.61
$ws.Range("M131").Value = -107376012
$ws.Range("N131").Value = -219233.61
$ws.Range("H133").Value = 5627.8335
$ws.Range("I133").Value = 4153.4
$ws.Range("J133").Value = 13000
$ws.Range("K133").Value = 12460.2
$ws.Range("L133").Value = 39000
$ws.Range("M133").Value = -7400.199999999999
$ws.Range("N133").Value = -49120
$ws.Range("H134").Value = 2565.2
$ws.Range("I134").Value = 2565.2
$ws.Range("K134").Value = 7695.599999999999
$ws.Range("M134").Value = -2625.599999999999
$ws.Range("H138").Value = 767.6667
$ws.Range("I138").Value = 767.6667
$ws.Range("K138").Value = 2303.0001
$ws.Range("M138").Value = 2836.9999
$ws.Range("H139").Value = 64052.438
$ws.Range("I139").Value = 67789.266
$ws.Range("K139").Value = 203367.798
$ws.Range("M139").Value = -198227.798

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 370.94116
$ws.Range("I97").Value = 361.35715
$ws.Range("K97").Value = 361.35715
$ws.Range("M97").Value = 134.64285
$ws.Range("H102").Value = 1739.2858
$ws.Range("I102").Value = 1716.3684
$ws.Range("J102").Value = 1957
$ws.Range("K102").Value = 1716.3684
$ws.Range("L102").Value = 1957
$ws.Range("M102").Value = -94.36840000000007
$ws.Range("N102").Value = -5201
$ws.Range("H107").Value = 2193.3
$ws.Range("I107").Value = 1836.4546
$ws.Range("J107").Value = 2399.8948
$ws.Range("K107").Value = 1836.4546
$ws.Range("L107").Value = 2399.8948
$ws.Range("M107").Value = 83.54539999999997
$ws.Range("N107").Value = -6239.8948
$ws.Range("H126").Value = 61999.8
$ws.Range("I126").Value = 1999.6666
$ws.Range("J126").Value = 152000
$ws.Range("K126").Value = 5998.9998
$ws.Range("L126").Value = 456000
$ws.Range("M126").Value = -3528.9998
$ws.Range("N126").Value = -460940
$ws.Range("H132").Value = 8348.5625
$ws.Range("I132").Value = 9198.308000000001
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 27594.924
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -25064.924
$ws.Range("N132").Value = -19059.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1124.4736
$ws.Range("I16").Value = 1150.9333
$ws.Range("K16").Value = 1150.9333
$ws.Range("M16").Value = -980.9332999999999
$ws.Range("H22").Value = 1342.7142
$ws.Range("I22").Value = 977.6667
$ws.Range("J22").Value = 1999.8
$ws.Range("K22").Value = 977.6667
$ws.Range("L22").Value = 1999.8
$ws.Range("M22").Value = -682.6667
$ws.Range("N22").Value = -2589.8
$ws.Range("H27").Value = 1342.7142
$ws.Range("I27").Value = 977.6667
$ws.Range("J27").Value = 1999.8
$ws.Range("K27").Value = 977.6667
$ws.Range("L27").Value = 1999.8
$ws.Range("M27").Value = -870.6667
$ws.Range("N27").Value = -2213.8
$ws.Range("H40").Value = 6853
$ws.Range("I40").Value = 6710.25
$ws.Range("K40").Value = 6710.25
$ws.Range("M40").Value = -6574.25
$ws.Range("H46").Value = 1164
$ws.Range("I46").Value = 1164
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1164
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -976
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 363.13635
$ws.Range("J55").Value = 467.625
$ws.Range("L55").Value = 467.625
$ws.Range("N55").Value = -813.625
$ws.Range("H61").Value = 37044360
$ws.Range("I61").Value = 41673904
$ws.Range("K61").Value = 41673904
$ws.Range("M61").Value = -41673702
$ws.Range("H68").Value = 2270.1
$ws.Range("I68").Value = 2014.5714
$ws.Range("J68").Value = 2866.3333
$ws.Range("K68").Value = 2014.5714
$ws.Range("L68").Value = 2866.3333
$ws.Range("M68").Value = -1265.5714
$ws.Range("N68").Value = -4364.3333
$ws.Range("H71").Value = 2270.1
$ws.Range("I71").Value = 2014.5714
$ws.Range("J71").Value = 2866.3333
$ws.Range("K71").Value = 10072.857
$ws.Range("L71").Value = 14331.6665
$ws.Range("M71").Value = -6328.857
$ws.Range("N71").Value = -21819.6665
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 37044360
$ws.Range("I113").Value = 41673904
$ws.Range("K113").Value = 41673904
$ws.Range("M113").Value = -41671734
$ws.Range("H122").Value = 4321.778
$ws.Range("I122").Value = 4099.273
$ws.Range("J122").Value = 4671.4287
$ws.Range("K122").Value = 12297.819
$ws.Range("L122").Value = 14014.2861
$ws.Range("M122").Value = -9847.819
$ws.Range("N122").Value = -18914.2861
$ws.Range("H132").Value = 4387.25
$ws.Range("I132").Value = 2882.0908
$ws.Range("J132").Value = 7698.6
$ws.Range("K132").Value = 8646.2724
$ws.Range("L132").Value = 23095.8
$ws.Range("M132").Value = -6116.2724
$ws.Range("N132").Value = -28155.8
$ws.Range("H136").Value = 3874.5
$ws.Range("J136").Value = 4999.8335
$ws.Range("L136").Value = 14999.5005
$ws.Range("N136").Value = -20099.5005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8822.538
$ws.Range("J81").Value = 15396.286
$ws.Range("L81").Value = 30792.572
$ws.Range("N81").Value = -32914.572
$ws.Range("H84").Value = 8822.538
$ws.Range("J84").Value = 15396.286
$ws.Range("L84").Value = 153962.86
$ws.Range("N84").Value = -164570.86
$ws.Range("H94").Value = 26386.666
$ws.Range("J94").Value = 26386.666
$ws.Range("L94").Value = 26386.666
$ws.Range("N94").Value = -28188.666
$ws.Range("H113").Value = 1946.9166
$ws.Range("I113").Value = 1542.75
$ws.Range("K113").Value = 4628.25
$ws.Range("M113").Value = -2458.25
$ws.Range("H122").Value = 1996.5588
$ws.Range("I122").Value = 1563.421
$ws.Range("K122").Value = 4690.263
$ws.Range("M122").Value = -2240.263
$ws.Range("H126").Value = 2830.0908
$ws.Range("I126").Value = 1713.3
$ws.Range("K126").Value = 5139.9
$ws.Range("M126").Value = -2669.9
$ws.Range("H132").Value = 4850.186
$ws.Range("I132").Value = 5143.7026
$ws.Range("J132").Value = 3040.1667
$ws.Range("K132").Value = 15431.1078
$ws.Range("L132").Value = 9120.500100000001
$ws.Range("M132").Value = -12901.1078
$ws.Range("N132").Value = -14180.5001

